$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Remove the "Фтор" / "Fluorine" / "FLD" nutrient row entirely (original row 17).
$ws.Range("A17").EntireRow.Delete()

# 2. Remove the "Ситість" / "Fullness" nutrient row entirely.
#    After the deletion above, everything shifted up by one row, so the row that
#    used to be 34 is now 33.
$ws.Range("A33").EntireRow.Delete()

# 3. Fill in the previously-missing "mda" (column E) values for several nutrients.
$ws.Range("E6").Value = 1000000
$ws.Range("E7").Value = 1000000
$ws.Range("E8").Value = 1000000
$ws.Range("E12").Value = 1000000
$ws.Range("E13").Value = 1000000
$ws.Range("E24").Value = 1000000000
$ws.Range("E26").Value = 1000000
$ws.Range("E27").Value = 1000000
$ws.Range("E29").Value = 1000000
$ws.Range("E30").Value = 1000000
